# Updated symbol list on Mon Dec 19 10:16:40 UTC 2022 with GitHub Actions
#
# Applies the per-row "Price" (column D) and "Hora" (column G) refresh
# pulled from coinranking.com, plus a reordering of two rows (BKEXToken
# and CEJI swapped position in the source feed).

function Set-TextValue($ws, $cell, $val) {
    # Force the cell to stay a Text-typed value (the sheet stores every
    # numeric-looking field, like Price/Hora, as literal text) and reset
    # the style back to Normal so no stray number-format gets attached.
    $ws.Range($cell).NumberFormat = "@"
    $ws.Range($cell).Value = $val
    $ws.Range($cell).Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 42 / Row 43 swapped places in the source feed (BKEXToken <-> CEJI) ---
Set-TextValue $ws "B42" "CEJI"
Set-TextValue $ws "C42" "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue $ws "D42" "0.003501"
Set-TextValue $ws "E42" "41CEJICEJI"

Set-TextValue $ws "B43" "BKEXToken"
Set-TextValue $ws "C43" "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue $ws "D43" "0.1044"
Set-TextValue $ws "E43" "42BKEXTokenBKK"

# --- Price (D) and Hora (G) refresh for every data row ---
Set-TextValue $ws "D2" "247.61"
Set-TextValue $ws "G2" "10"
Set-TextValue $ws "G3" "10"
Set-TextValue $ws "D4" "5.479"
Set-TextValue $ws "G4" "10"
Set-TextValue $ws "D5" "0.05692"
Set-TextValue $ws "G5" "10"
Set-TextValue $ws "D6" "3.380"
Set-TextValue $ws "G6" "10"
Set-TextValue $ws "D7" "0.8055"
Set-TextValue $ws "G7" "10"
Set-TextValue $ws "D8" "1.037"
Set-TextValue $ws "G8" "10"
Set-TextValue $ws "D9" "0.1522"
Set-TextValue $ws "G9" "10"
Set-TextValue $ws "D10" "0.07317"
Set-TextValue $ws "G10" "10"
Set-TextValue $ws "D11" "0.03165"
Set-TextValue $ws "G11" "10"
Set-TextValue $ws "D12" "0.02967"
Set-TextValue $ws "G12" "10"
Set-TextValue $ws "D13" "0.09295"
Set-TextValue $ws "G13" "10"
Set-TextValue $ws "D14" "3.432"
Set-TextValue $ws "G14" "10"
Set-TextValue $ws "D15" "0.001646"
Set-TextValue $ws "G15" "10"
Set-TextValue $ws "D16" "0.04725"
Set-TextValue $ws "G16" "10"
Set-TextValue $ws "D17" "0.0005871"
Set-TextValue $ws "G17" "10"
Set-TextValue $ws "D18" "0.006357"
Set-TextValue $ws "G18" "10"
Set-TextValue $ws "D19" "0.005042"
Set-TextValue $ws "G19" "10"
Set-TextValue $ws "D20" "0.001046"
Set-TextValue $ws "G20" "10"
Set-TextValue $ws "D21" "0.0001500"
Set-TextValue $ws "G21" "10"
Set-TextValue $ws "D22" "0.0003136"
Set-TextValue $ws "G22" "10"
Set-TextValue $ws "D23" "3.772"
Set-TextValue $ws "G23" "10"
Set-TextValue $ws "D24" "6.431"
Set-TextValue $ws "G24" "10"
Set-TextValue $ws "D25" "2.112"
Set-TextValue $ws "G25" "10"
Set-TextValue $ws "D26" "0.3280"
Set-TextValue $ws "G26" "10"
Set-TextValue $ws "D27" "0.1300"
Set-TextValue $ws "G27" "10"
Set-TextValue $ws "G28" "10"
Set-TextValue $ws "G29" "10"
Set-TextValue $ws "G30" "10"
Set-TextValue $ws "G31" "10"
Set-TextValue $ws "G32" "10"
Set-TextValue $ws "G33" "10"
Set-TextValue $ws "G34" "10"
Set-TextValue $ws "G35" "10"
Set-TextValue $ws "G36" "10"
Set-TextValue $ws "G37" "10"
Set-TextValue $ws "G38" "10"
Set-TextValue $ws "G39" "10"
Set-TextValue $ws "D40" "0.04107"
Set-TextValue $ws "G40" "10"
Set-TextValue $ws "D41" "0.006932"
Set-TextValue $ws "G41" "10"
Set-TextValue $ws "G42" "10"
Set-TextValue $ws "G43" "10"
Set-TextValue $ws "D44" "0.008634"
Set-TextValue $ws "G44" "10"
Set-TextValue $ws "D45" "0.00005827"
Set-TextValue $ws "G45" "10"
Set-TextValue $ws "G46" "10"
Set-TextValue $ws "D47" "0.0005501"
Set-TextValue $ws "G47" "10"
Set-TextValue $ws "D48" "0.6826"
Set-TextValue $ws "G48" "10"
Set-TextValue $ws "D49" "0.009428"
Set-TextValue $ws "G49" "10"
Set-TextValue $ws "D50" "0.00002101"
Set-TextValue $ws "G50" "10"
Set-TextValue $ws "G51" "10"
